$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.865.56'
$ws.Range("E2").Value = '  +2.47%  '

$ws.Range("D3").Value = '2.563.22'
$ws.Range("E3").Value = '  +1.84%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Value = '601.15'
$ws.Range("E5").Value = '  +1.92%  '

$ws.Range("D6").Value = '178.69'
$ws.Range("E6").Value = '  +0.60%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("E8").Value = '  +0.53%  '

$ws.Range("D9").Value = '2.562.42'
$ws.Range("E9").Value = '  +1.79%  '

$ws.Range("D10").Value = '0.159'
$ws.Range("E10").Value = '  +11.57%  '

$ws.Range("D12").Value = '0.345'
$ws.Range("E12").Value = '  +1.24%  '

$ws.Range("D13").Value = '5.02'
$ws.Range("E13").Value = '  +1.34%  '

$ws.Range("D14").Value = '3.042.40'
$ws.Range("E14").Value = '  +5.77%  '

$ws.Range("D15").Value = '26.34'
$ws.Range("E15").Value = '  +1.98%  '

$ws.Range("E16").Value = '  +5.22%  '

$ws.Range("D17").Value = '69.812.07'
$ws.Range("E17").Value = '  +2.52%  '

$ws.Range("D18").Value = '2.568.66'
$ws.Range("E18").Value = '  +1.26%  '

$ws.Range("D19").Value = '7.69'
$ws.Range("E19").Value = '  +1.88%  '

$ws.Range("D20").Value = '11.17'
$ws.Range("E20").Value = '  +1.19%  '

$ws.Range("D21").Value = '365.11'
$ws.Range("E21").Value = '  +3.55%  '

$ws.Range("E22").Value = '  +2.10%  '

$ws.Range("E23").Value = '  -0.19%  '

$ws.Range("D24").Value = '70.74'
$ws.Range("E24").Value = '  +0.06%  '

$ws.Range("D25").Value = '4.31'

$ws.Range("E26").Value = '  -1.49%  '

$ws.Range("D27").Value = '9.26'
$ws.Range("E27").Value = '  +1.12%  '

$ws.Range("D28").Value = '2.695.08'
$ws.Range("E28").Value = '  +2.06%  '

$ws.Range("D29").Value = '0.996'
$ws.Range("E29").Value = '  +0.56%  '

$ws.Range("D30").Value = '0.0₃0931'
$ws.Range("E30").Value = '  +1.15%  '

$ws.Range("D31").Value = '520.19'
$ws.Range("E31").Value = '  +1.89%  '

$ws.Range("D32").Value = '7.79'
$ws.Range("E32").Value = '  -1.43%  '

$ws.Range("E33").Value = '  -0.11%  '

$ws.Range("E34").Value = '  +1.38%  '

$ws.Range("E35").Value = '  +0.03%  '

$ws.Range("D36").Value = '0.120'
$ws.Range("E36").Value = '  -1.52%  '

$ws.Range("D37").Value = '163.08'
$ws.Range("E37").Value = '  -1.17%  '

$ws.Range("D38").Value = '18.99'
$ws.Range("E38").Value = '  +2.94%  '

$ws.Range("D39").Value = '18.93'
$ws.Range("E39").Value = '  +1.41%  '

$ws.Range("E40").Value = '  +0.56%  '

$ws.Range("E41").Value = '  +0.91%  '

$ws.Range("E42").Value = '  -0.01%  '

$ws.Range("D43").Value = '4.94'
$ws.Range("E43").Value = '  +0.69%  '

$ws.Range("D44").Value = '0.326'
$ws.Range("E44").Value = '  -1.27%  '

$ws.Range("D45").Value = '2.47'
$ws.Range("E45").Value = '  -0.69%  '

$ws.Range("D46").Value = '39.06'
$ws.Range("E46").Value = '  +0.33%  '

$ws.Range("D47").Value = '153.08'
$ws.Range("E47").Value = '  +4.12%  '

$ws.Range("D48").Value = '3.63'
$ws.Range("E48").Value = '  +2.17%  '

$ws.Range("D49").Value = '0.524'
$ws.Range("E49").Value = '  +0.51%  '

$ws.Range("D50").Value = '0.0₆0259'
$ws.Range("E50").Value = '  -0.68%  '

$ws.Range("E51").Value = '  +1.56%  '
